$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

# Row 100
$ws.Range("H100").Value = 1246.1818
$ws.Range("I100").Value = 1673.75
$ws.Range("J100").Value = 1001.8571
$ws.Range("K100").Value = 1673.75
$ws.Range("L100").Value = 1001.8571
$ws.Range("M100").Value = -1132.75
$ws.Range("N100").Value = -2083.8571

# Row 113
$ws.Range("H113").Value = 61153.234
$ws.Range("I113").Value = 169734.17
$ws.Range("J113").Value = 1927.2727
$ws.Range("K113").Value = 169734.17
$ws.Range("L113").Value = 1927.2727
$ws.Range("M113").Value = -166480.17
$ws.Range("N113").Value = -8435.2727

# Row 129
$ws.Range("H129").Value = 305016.28
$ws.Range("I129").Value = 10444.9
$ws.Range("J129").Value = 386841.66
$ws.Range("K129").Value = 31334.7
$ws.Range("L129").Value = 1160524.98
$ws.Range("M129").Value = -26334.7
$ws.Range("N129").Value = -1170524.98

# Row 137
$ws.Range("H137").Value = 1482.2122
$ws.Range("I137").Value = 1212.1538
$ws.Range("K137").Value = 3636.4614
$ws.Range("M137").Value = -1086.4614

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 75
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("M4").Value = 66

# Row 5
$ws.Range("H5").Value = 333393.34
$ws.Range("I5").Value = 1000000
$ws.Range("K5").Value = 1000000
$ws.Range("M5").Value = -999888

# Row 32
$ws.Range("H32").Value = 23156.012
$ws.Range("I32").Value = 3930.2207
$ws.Range("K32").Value = 3930.2207
$ws.Range("M32").Value = -3643.2207

# Row 130
$ws.Range("H130").Value = 48985
$ws.Range("J130").Value = 48985
$ws.Range("L130").Value = 48985
$ws.Range("N130").Value = -59025

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 333393.34
$ws.Range("I4").Value = 1000000
$ws.Range("K4").Value = 1000000
$ws.Range("M4").Value = -999885

# Row 86
$ws.Range("H86").Value = 66631.84
$ws.Range("I86").Value = 95493.69500000001
$ws.Range("J86").Value = 4097.8335
$ws.Range("K86").Value = 95493.69500000001
$ws.Range("L86").Value = 4097.8335
$ws.Range("M86").Value = -94370.69500000001
$ws.Range("N86").Value = -6343.8335

# Row 89
$ws.Range("H89").Value = 66631.84
$ws.Range("I89").Value = 95493.69500000001
$ws.Range("J89").Value = 4097.8335
$ws.Range("K89").Value = 477468.475
$ws.Range("L89").Value = 20489.1675
$ws.Range("M89").Value = -471852.475
$ws.Range("N89").Value = -31721.1675

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 23941.338
$ws.Range("J31").Value = 37074.043
$ws.Range("L31").Value = 37074.043
$ws.Range("N31").Value = -37664.043

# Row 34
$ws.Range("H34").Value = 23941.338
$ws.Range("J34").Value = 37074.043
$ws.Range("L34").Value = 37074.043
$ws.Range("N34").Value = -37478.043

# Row 58
$ws.Range("H58").Value = 6083.3125
$ws.Range("I58").Value = 1167.5
$ws.Range("J58").Value = 20830.75
$ws.Range("K58").Value = 1167.5
$ws.Range("L58").Value = 20830.75
$ws.Range("M58").Value = -964.5
$ws.Range("N58").Value = -21236.75

# Row 130
$ws.Range("H130").Value = 47986.332
$ws.Range("J130").Value = 47986.332
$ws.Range("L130").Value = 47986.332
$ws.Range("N130").Value = -58026.332

# Row 136
$ws.Range("H136").Value = 6083.3125
$ws.Range("I136").Value = 1167.5
$ws.Range("J136").Value = 20830.75
$ws.Range("K136").Value = 3502.5
$ws.Range("L136").Value = 62492.25
$ws.Range("M136").Value = -952.5
$ws.Range("N136").Value = -67592.25

$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Range("H62").Value = 2726.6
$ws.Range("I62").Value = 799.75
$ws.Range("J62").Value = 3427.2727
$ws.Range("K62").Value = 2399.25
$ws.Range("L62").Value = 10281.8181
$ws.Range("M62").Value = -1713.25
$ws.Range("N62").Value = -11653.8181

# Row 63
$ws.Range("J63").Value = 1200
$ws.Range("L63").Value = 3600
$ws.Range("N63").Value = -5098

# Row 64
$ws.Range("H64").Value = 2491.6667
$ws.Range("I64").Value = 2250
$ws.Range("J64").Value = 2560.7144
$ws.Range("K64").Value = 6750
$ws.Range("L64").Value = 7682.1432
$ws.Range("M64").Value = -6480
$ws.Range("N64").Value = -8222.143199999999

# Row 65
$ws.Range("H65").Value = 2726.6
$ws.Range("I65").Value = 799.75
$ws.Range("J65").Value = 3427.2727
$ws.Range("K65").Value = 7197.75
$ws.Range("L65").Value = 30845.4543
$ws.Range("M65").Value = -3765.75
$ws.Range("N65").Value = -37709.4543

# Row 66
$ws.Range("J66").Value = 1200
$ws.Range("L66").Value = 10800
$ws.Range("N66").Value = -18288

# Row 67
$ws.Range("H67").Value = 2491.6667
$ws.Range("I67").Value = 2250
$ws.Range("J67").Value = 2560.7144
$ws.Range("K67").Value = 6750
$ws.Range("L67").Value = 7682.1432
$ws.Range("M67").Value = -5814
$ws.Range("N67").Value = -9554.143199999999

# Row 70
$ws.Range("H70").Value = 127176.5
$ws.Range("J70").Value = 2600
$ws.Range("L70").Value = 7800
$ws.Range("N70").Value = -8430

# Row 73
$ws.Range("H73").Value = 127176.5
$ws.Range("J73").Value = 2600
$ws.Range("L73").Value = 7800
$ws.Range("N73").Value = -9984

# Row 74
$ws.Range("H74").Value = 4000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 12000
$ws.Range("N74").Value = -14122
$ws.Range("M74").ClearContents()

# Row 77
$ws.Range("H77").Value = 4000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 36000
$ws.Range("N77").Value = -46608
$ws.Range("M77").ClearContents()

# Row 119
$ws.Range("H119").Value = 68566.87
$ws.Range("I119").Value = 100537.5
$ws.Range("J119").Value = 4625.6
$ws.Range("K119").Value = 301612.5
$ws.Range("L119").Value = 13876.8
$ws.Range("M119").Value = -296774.5
$ws.Range("N119").Value = -23552.8

# Row 131
$ws.Range("H131").Value = 813.45
$ws.Range("J131").Value = 848.8461
$ws.Range("L131").Value = 2546.5383
$ws.Range("N131").Value = -12626.5383

$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 11752
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 11752
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 11752
$ws.Range("N52").Value = -12270
$ws.Range("M52").ClearContents()

# Row 70
$ws.Range("H70").Value = 83441.92
$ws.Range("I70").Value = 117138.11
$ws.Range("J70").Value = 7625.5
$ws.Range("K70").Value = 117138.11
$ws.Range("L70").Value = 7625.5
$ws.Range("M70").Value = -116868.11
$ws.Range("N70").Value = -8165.5

# Row 73
$ws.Range("H73").Value = 83441.92
$ws.Range("I73").Value = 117138.11
$ws.Range("J73").Value = 7625.5
$ws.Range("K73").Value = 117138.11
$ws.Range("L73").Value = 7625.5
$ws.Range("M73").Value = -116202.11
$ws.Range("N73").Value = -9497.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 779150.7
$ws.Range("I46").Value = 369.5
$ws.Range("J46").Value = 920747.25
$ws.Range("K46").Value = 369.5
$ws.Range("L46").Value = 369.5
$ws.Range("M46").Value = -181.5
$ws.Range("N46").Value = -921123.25

# Row 68
$ws.Range("H68").Value = 3059.2
$ws.Range("I68").Value = 1516.5
$ws.Range("J68").Value = 4087.6667
$ws.Range("K68").Value = 1516.5
$ws.Range("L68").Value = 4087.6667
$ws.Range("M68").Value = -767.5
$ws.Range("N68").Value = -5585.6667

# Row 71
$ws.Range("H71").Value = 3059.2
$ws.Range("I71").Value = 1516.5
$ws.Range("J71").Value = 4087.6667
$ws.Range("K71").Value = 7582.5
$ws.Range("L71").Value = 20438.3335
$ws.Range("M71").Value = -3838.5
$ws.Range("N71").Value = -27926.3335

# Row 82
$ws.Range("H82").Value = 1213.7
$ws.Range("I82").Value = 804.44446
$ws.Range("J82").Value = 1548.5454
$ws.Range("K82").Value = 804.44446
$ws.Range("L82").Value = 1548.5454
$ws.Range("M82").Value = -443.44446
$ws.Range("N82").Value = -2270.5454

# Row 85
$ws.Range("H85").Value = 1213.7
$ws.Range("I85").Value = 804.44446
$ws.Range("J85").Value = 1548.5454
$ws.Range("K85").Value = 804.44446
$ws.Range("L85").Value = 1548.5454
$ws.Range("M85").Value = 443.55554
$ws.Range("N85").Value = -4044.5454
